# Apply hybrid bold + color (#2C3E50) highlighting to quantitative metrics
# (percentages, dollar amounts, large numbers) inside specific bullet
# paragraphs of the resume, per the commit:
# "Implement quantitative metrics highlighting across all resume formats"

$d = $word.ActiveDocument

# Word VT_COLOR decimal for hex 2C3E50 (R=0x2C,G=0x3E,B=0x50 -> R + G*256 + B*65536)
$highlightColor = 5258796

function Highlight-Metric {
    param(
        [int]$paraIndex,
        [string]$metricText
    )
    $p = $d.Paragraphs.Item($paraIndex)
    $scopeStart = $p.Range.Start
    $scopeEnd = $p.Range.End
    $rng = $d.Range($scopeStart, $scopeEnd)
    $found = $rng.Find.Execute($metricText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Font.Bold = 1
        $rng.Font.Color = $highlightColor
    }
}

# Paragraph 9: "...classification accuracy from 23% to 64%"
Highlight-Metric 9 '23%'
Highlight-Metric 9 '64%'

# Paragraph 11: "Achieved 87% ... standard of 71%, ... margins from ±4.2% to ±2.1%"
Highlight-Metric 11 '87%'
Highlight-Metric 11 '71%'
Highlight-Metric 11 '±4.2%'
Highlight-Metric 11 '±2.1%'

# Paragraph 31: "...bids from 1,200 vendors..."
Highlight-Metric 31 '1,200'

# Paragraph 46: "...became the $400M Polling Consortium Database ... valued at $1B+"
Highlight-Metric 46 '$400M'
Highlight-Metric 46 '$1B'

# Paragraph 63: "Algorithm reduced mapping costs by 73.5%, saving ... $4.7M"
Highlight-Metric 63 '73.5%'
Highlight-Metric 63 '$4.7M'

# Paragraph 65: "Achieved 87% prediction accuracy ... standard of 71%"
Highlight-Metric 65 '87%'
Highlight-Metric 65 '71%'
